$d = $word.ActiveDocument

# Find the existing author name "Nakul Mody" so we can append the new
# co-author right after it, inside the same paragraph / run of text.
$range = $d.Content
$found = $range.Find.Execute("Nakul Mody", $true, $false, $false, $false,
                              $false, $true, 1, $false, "", 0)

# $range now covers "Nakul Mody" after a successful Find.Execute; collapse
# to its end point so the insertion lands immediately after the name.
$range.Collapse(0)  # wdCollapseEnd = 0

# Insert the new co-author text.
$range.InsertAfter(", Jesse Oh")

# Match the run formatting used by the rest of the author line
# (Times New Roman ascii/eastAsia/hAnsi/cs, bold, size 28 half-points -> 14pt, kern 0).
$range.Font.Name = "Times New Roman"
$range.Font.NameAscii = "Times New Roman"
$range.Font.NameFarEast = "Times New Roman"
$range.Font.NameBi = "Times New Roman"
$range.Font.NameOther = "Times New Roman"
$range.Font.Bold = $true
$range.Font.Size = 14
$range.Font.Kerning = 0
